$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "GossA-HW30.xpc" to "GossA"
$ws.Name = "GossA"

# Copy formatting (bold font, border, centered alignment) from A15 to the new A16 cell
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats

# Add new row 16 data (14th additional Gaussian-quadrature scheme result row)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.038396101144211
$ws.Range("D16").Value = 0.9568473323611687
$ws.Range("E16").Value = 0.9990556871472746
$ws.Range("F16").Value = 0.9826358984138658
$ws.Range("G16").Value = 1.038396101144211
$ws.Range("H16").Value = 0.9568473323611687
$ws.Range("I16").Value = 1.007458715600644
$ws.Range("J16").Value = 0.9807765785170033
$ws.Range("K16").Value = 1.009811656777114
$ws.Range("L16").Value = 0.9676101139399597
$ws.Range("M16").Value = 1.038396101144211
$ws.Range("N16").Value = 0.9779515097542216
$ws.Range("O16").Value = 0.99423375476663
$ws.Range("P16").Value = 0.9928240104876551
